$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Page URL"
$ws.Range("B1").Value = "Test Case"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Comments"
